# Generate Report for Handback
# For the "cbae1ca0-eb85-4205-a5eb-958f2e3718c3" entry (row 5 in both the
# zh-cn and de-de sheets) record that the handback that was received is
# stale: populate the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns, widen those
# columns, and turn the "Latest Target File" cell into a hyperlink back to
# the handback markdown file (same target as the existing A5 hyperlink).

$wb = $excel.ActiveWorkbook

$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0203d389af8492a7798bf300fa970b8c35f2eeb/e2e/cbae1ca0-eb85-4205-a5eb-958f2e3718c3.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7464524c8766c28c6f310edab1ece23869c33d7d/e2e/cbae1ca0-eb85-4205-a5eb-958f2e3718c3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0203d389af8492a7798bf300fa970b8c35f2eeb/e2e/cbae1ca0-eb85-4205-a5eb-958f2e3718c3.md."

function Update-LanguageSheet($sheetName, $handbackFileName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen Latest Target File (I), Latest Handback File (J) and Error Detail (P)
    $ws.Columns.Item(9).ColumnWidth = 39.14
    $ws.Columns.Item(10).ColumnWidth = 39.14
    $ws.Columns.Item(16).ColumnWidth = 39.14

    # Latest Handback File (J5)
    $ws.Range("J5").Value = $handbackFileName

    # Latest Handback DateTime (K5) - stored as text, matching the sheet's convention
    $ws.Range("K5").Value = $handbackDateTime

    # Error Detail (P5)
    $ws.Range("P5").Value = $errorDetail

    # Latest Target File (I5) - link back to the handback markdown, same
    # target as the existing A5 hyperlink for this row.
    $ws.Hyperlinks.Add($ws.Range("I5"), $latestHandbackUrl, $null, $null, "cbae1ca0-eb85-4205-a5eb-958f2e3718c3.md")
}

Update-LanguageSheet "zh-cn" "cbae1ca0-eb85-4205-a5eb-958f2e3718c3.b3650fce2b8ee257805d1e37e3c4059308a06ded.zh-cn.xlf" "2016-08-31 07:21:06"
Update-LanguageSheet "de-de" "cbae1ca0-eb85-4205-a5eb-958f2e3718c3.b3650fce2b8ee257805d1e37e3c4059308a06ded.de-de.xlf" "2016-08-31 07:21:27"
